$d = $word.ActiveDocument

$d.Content.Find.Execute("217÷8=27, 1", $true, $false, $false, $false, $false, $true, 1, $false, "182÷5=36, 2", 2)
$d.Content.Find.Execute("684÷5=136, 4", $true, $false, $false, $false, $false, $true, 1, $false, "434÷6=72, 2", 2)
$d.Content.Find.Execute("433÷3=144, 1", $true, $false, $false, $false, $false, $true, 1, $false, "193÷9=21, 4", 2)
$d.Content.Find.Execute("780÷8=97, 4", $true, $false, $false, $false, $false, $true, 1, $false, "338÷5=67, 3", 2)
$d.Content.Find.Execute("875÷4=218, 3", $true, $false, $false, $false, $false, $true, 1, $false, "741÷8=92, 5", 2)
$d.Content.Find.Execute("874÷2=437, 0", $true, $false, $false, $false, $false, $true, 1, $false, "825÷6=137, 3", 2)
$d.Content.Find.Execute("877÷4=219, 1", $true, $false, $false, $false, $false, $true, 1, $false, "228÷6=38, 0", 2)
$d.Content.Find.Execute("998÷7=142, 4", $true, $false, $false, $false, $false, $true, 1, $false, "523÷8=65, 3", 2)
$d.Content.Find.Execute("662÷8=82, 6", $true, $false, $false, $false, $false, $true, 1, $false, "898÷4=224, 2", 2)
$d.Content.Find.Execute("233÷9=25, 8", $true, $false, $false, $false, $false, $true, 1, $false, "390÷3=130, 0", 2)
$d.Content.Find.Execute("500÷4=125, 0", $true, $false, $false, $false, $false, $true, 1, $false, "526÷2=263, 0", 2)
$d.Content.Find.Execute("818÷2=409, 0", $true, $false, $false, $false, $false, $true, 1, $false, "230÷6=38, 2", 2)
$d.Content.Find.Execute("694÷5=138, 4", $true, $false, $false, $false, $false, $true, 1, $false, "218÷3=72, 2", 2)
$d.Content.Find.Execute("737÷2=368, 1", $true, $false, $false, $false, $false, $true, 1, $false, "609÷8=76, 1", 2)
$d.Content.Find.Execute("524÷7=74, 6", $true, $false, $false, $false, $false, $true, 1, $false, "830÷6=138, 2", 2)
$d.Content.Find.Execute("612÷5=122, 2", $true, $false, $false, $false, $false, $true, 1, $false, "839÷4=209, 3", 2)
$d.Content.Find.Execute("856÷6=142, 4", $true, $false, $false, $false, $false, $true, 1, $false, "549÷8=68, 5", 2)
$d.Content.Find.Execute("236÷8=29, 4", $true, $false, $false, $false, $false, $true, 1, $false, "171÷3=57, 0", 2)
$d.Content.Find.Execute("756÷3=252, 0", $true, $false, $false, $false, $false, $true, 1, $false, "624÷6=104, 0", 2)
$d.Content.Find.Execute("861÷4=215, 1", $true, $false, $false, $false, $false, $true, 1, $false, "593÷6=98, 5", 2)
$d.Content.Find.Execute("142÷6=23, 4", $true, $false, $false, $false, $false, $true, 1, $false, "764÷2=382, 0", 2)
$d.Content.Find.Execute("511÷5=102, 1", $true, $false, $false, $false, $false, $true, 1, $false, "816÷2=408, 0", 2)
$d.Content.Find.Execute("250÷7=35, 5", $true, $false, $false, $false, $false, $true, 1, $false, "506÷7=72, 2", 2)
$d.Content.Find.Execute("846÷6=141, 0", $true, $false, $false, $false, $false, $true, 1, $false, "827÷5=165, 2", 2)
$d.Content.Find.Execute("102÷5=20, 2", $true, $false, $false, $false, $false, $true, 1, $false, "102÷7=14, 4", 2)
